# Tuesday 16 July signoff backup
# Apply the "Daily report" sheet updates (new row-4 cells, new row 5, taller
# rows) and refresh the saved selections on both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Daily report" sheet
# ---------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily report")
$daily.Activate()

# --- Row 4: add entries in E4/F4, and give C4 the same wrap+vcenter
#     style already used by B4/D4.
$c4 = $daily.Range("C4")
$c4.WrapText = $true
$c4.VerticalAlignment = -4108   # xlCenter

$e4 = $daily.Range("E4")
$e4.Value = @'
Tested and confirmed GUI.py and Unity project work, focused on Steam Audio quirks, decided to create a fresh Steam Audio project for comparison, noted exact dimension reproduction for monodepth isn't crucial, explored enhance360.py arguments to strengthen depth maps without complex preprocessing, found increasing baseline value helps retain detail, faced issues at baseline 5.264, followed digital audio fundamentals on YouTube, planning to use Audacity for waveform comparison, considered simplifying preprocessing by removing top and bottom modifications only, realized mirror fixes were incidental, commented out plane detection in enhance360.py for better results, but still need to address some depth loss.
'@

$f4 = $daily.Range("F4")
$f4.Value = @'
Worked from home, cleaned up the GitHub repo for better pulls, faced issues with Unity project .meta files and prefab initialization, ran into Git LFS bandwidth limits due to repeated .obj file pulls, decided to check for .obj file existence before import to avoid errors, planned to use a different storage provider for large files, decided to update the manual with a basic checklist for new PC/system setups, considered reactivating boostingMonocularDepth in 360monodepth for better depth images without enhance360.py, updated TODO to include fixing Steam Audio quirks, optimizing monodepth with boostingMonocularDepth, and removing Docker clutter.
'@

# --- Row 5: new entries in B5/C5. B5 needs the wrap+vcenter style.
$b5 = $daily.Range("B5")
$b5.WrapText = $true
$b5.VerticalAlignment = -4108   # xlCenter
$b5.Value = @'
Met with Dr. Hansung Kim at 10:30 am, discussed improving Unity project realism, adding collision for non-VR player controller, implementing binaural audio, and checking RIR calculations. Mona will WFH next week; asked necessary questions and to send recording to Dr. Hansung for feedback. Updated  Fixed Dockerfile build error, confirmed Monodepth v1.0 has different Image IDs, refactored code, compared enhance360.py with new enhance360mono.py, found enhance360mono better for depth details but less accurate for room structure. Tested on UL, found MonodepthBoosting not as good as manual edge optimization, generated meshes showed artifacts, MR scene better but not ideal. Plan to generate all scenes and consult Dr. Hansung Kim, identified V2's white balance issue, will figure out monoenhance tomorrow, and may fallback to enhance360.py baseline 2.264.
'@

$c5 = $daily.Range("C5")
$c5.WrapText = $true
$c5.Value = @'
Reconfigured StartScene to include Steam Audio Probe Batch, added baked source and listener components, changed behavior to baked as per Mona’s recommendation, exported audio using Wavepad or Audacity, added checks and dynamic project directory, confirmed settings with Mona’s suggestions. Found baked settings clearer in directivity and reverb compared to real-time, note to export active scenes after every bake. Encountered bug with old project directory, fixed by reimporting Steam Audio. Baked probe only has one level when using uniform floor mode, requiring player controller with collision for realism. Radeon Rays mode caused OpenCL error, unsure about TrueNext in reflection mode. Baking time should be minimal if correctly configured. Experienced bug with occlusion and distance attenuation in freecam mode. Decided to have an in-game indicator for Steam Audio status. Shifted focus to enhance360mono, created visualizer script, identified fix limits as the issue, promising results with fix_limits_mono from Claude. Signed off for the day.
'@

# --- Row heights to fit the new/expanded content.
$daily.Rows.Item(4).RowHeight = 261
$daily.Rows.Item(5).RowHeight = 319

# --- Scroll / selection bookkeeping matching the saved view state.
$daily.Range("D5").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1

# ---------------------------------------------------------------------
# "10 Week" sheet - just the saved selection moved from D4 to D5.
# ---------------------------------------------------------------------
$tenWeek = $wb.Worksheets.Item("10 Week")
$tenWeek.Activate()
$tenWeek.Range("D5").Select()

# Restore "Daily report" as the active/visible tab.
$daily.Activate()
